$wb = $excel.ActiveWorkbook

# --- Sheet 1: rename Hoja1 -> Recursos, drop the "Descripcion Tarea" column (E) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Recursos"
$ws1.Columns("E:E").Delete()

# Leave the selection on the (now empty) column E, as in the authored file.
$ws1.Columns("E:E").Select() | Out-Null

# --- Sheet 2: brand new "Tareas" sheet, placed after Recursos ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Tareas"

$ws2.Range("A1").Value = "Categoria"
$ws2.Range("A1").Interior.ThemeColor = 3
$ws2.Range("B1").Value = "Tarea"
$ws2.Range("B1").Interior.ThemeColor = 3
$ws2.Range("C1").Value = "Código"
$ws2.Range("C1").Interior.ThemeColor = 3

$ws2.Columns("A:A").ColumnWidth = 31.736979166666668
$ws2.Columns("B:B").ColumnWidth = 30.022135416666668
$ws2.Columns("C:C").ColumnWidth = 11.022135416666666

$ws2.Range("A2:D363").Select()
